# "finish tests for computeEngine"
#
# On the "Compute Engine" checklist sheet, rows 2-14 (the "positive" cases)
# already have a note in column C saying the scenario was implemented in
# computeEngine.spec.ts. Rows 17-20 (the "negative" cases) were still
# missing that note. This finishes the checklist by adding the same note
# to those rows too, growing their row height to match the rest of the
# (now wrapped, taller) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compute Engine")
$ws.Select()

$note = "Реализовано в computeEngine.spec.ts"

# Copy the formatting already used for the "implemented" note (column C,
# e.g. C14) onto the four still-empty cells, then fill in the text.
$ws.Range("C14").Copy()
$ws.Range("C17:C20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C17").Value = $note
$ws.Range("C18").Value = $note
$ws.Range("C19").Value = $note
$ws.Range("C20").Value = $note

# With the extra wrapped column-C text, these rows now need the same
# taller row height already used elsewhere in the checklist.
$ws.Rows(17).RowHeight = 46
$ws.Rows(18).RowHeight = 46
$ws.Rows(19).RowHeight = 46
$ws.Rows(20).RowHeight = 46

$ws.Range("D12").Select()
